# The sheet contains one weekly record per row for "Espinaca" at the
# "Vega Modelo de Temuco" market. This commit adds a new (more recent)
# weekly observation. The new record is inserted as row 168 (pushing
# every existing row below it down by one), and is populated with the
# newest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 168; this shifts rows 168:205 down to 169:206
# (preserving their data/styles untouched) exactly like Excel's own
# Rows.Insert, and grows the sheet's used range to A1:R206 automatically.
$ws.Rows("168:168").Insert()

# Populate the newly inserted row 168 with the new weekly observation.
$ws.Range("A168").Value = 10
$ws.Range("B168").Value = "Vega Modelo de Temuco"
$ws.Range("C168").Value = "La Araucanía"
$ws.Range("D168").Value = 44855
$ws.Range("E168").Value = 9
$ws.Range("F168").Value = 100112012
$ws.Range("G168").Value = "Espinaca"
$ws.Range("H168").Value = "Sin especificar"
$ws.Range("I168").Value = "Primera"
$ws.Range("J168").Value = 30
$ws.Range("K168").Value = 10000
$ws.Range("L168").Value = 10000
$ws.Range("M168").Value = 10000
$ws.Range("N168").Value = "`$/docena de atados"
$ws.Range("O168").Value = "Región de La Araucanía"
$ws.Range("P168").Value = 3333
$ws.Range("Q168").Value = 3
$ws.Range("R168").Value = "Hortaliza"
